$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (id = 10) in row 11 - fill in the previously-empty row.
$ws.Range("B11").Value = 24
$ws.Range("C11").Formula = "=3*41548"
$ws.Range("D11").Formula = "=C11*F11"
$ws.Range("E11").Value = 484
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 1697
$ws.Range("H11").Value = 0.044398148148148145
$ws.Range("I11").Value = 6887
$ws.Range("J11").Value = "Vampiro"
$ws.Range("K11").Value = "Normal"
$ws.Range("L11").Value = 46014

# Match the number formats used by the rest of the table (copy formats only,
# values already set above so only formatting is pulled across).
$ws.Range("H10").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("L10").Copy()
$ws.Range("L11").PasteSpecial(-4122)

# Leave the selection where the author left it.
$null = $ws.Range("I11").Select()
